$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Defs" - rows 2-5 get refreshed _uid/_created/_updated values
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Defs")

$ws.Range("A2").Value = "lgs4l68f-wupv"
$ws.Range("B2").Value = "2023-04-22T15:18:42.639Z"
$ws.Range("C2").Value = "lgs4l690"

$ws.Range("A3").Value = "lgs4l68g-jlfh"
$ws.Range("B3").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C3").Value = "lgs4l68g"

$ws.Range("A4").Value = "lgs4l68g-095d"
$ws.Range("B4").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C4").Value = "lgs4l68g"

$ws.Range("A5").Value = "lgs4l690-0isn"
$ws.Range("B5").Value = "2023-04-22T15:18:42.639Z"
$ws.Range("C5").Value = "lgs4l690"

# ---------------------------------------------------------------------------
# Sheet "Point Defs" - drop the _format column (L), refresh uid/created/
# updated on existing rows, tweak a couple of values, mark row 4 deleted,
# and append row 5 as the replacement point def.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Point Defs")

$ws.Range("L1:L4").EntireColumn.Delete()

$ws.Range("A2").Value = "lgs4l68g-0qsw"
$ws.Range("B2").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C2").Value = "lgs4l68g"
$ws.Range("F2").Value = "ats6"

$ws.Range("A3").Value = "lgs4l68g-mj7s"
$ws.Range("B3").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C3").Value = "lgs4l68g"

$ws.Range("A4").Value = "lgs4l68g-hvoj"
$ws.Range("B4").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C4").Value = "lgs4l690"
$ws.Range("D4").Value = $true
$ws.Range("I4").Value = "Orig desc"

$ws.Range("A5").Value = "lgs4l690-9zv4e"
$ws.Range("B5").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C5").Value = "lgs4l690"
$ws.Range("D5").Value = $false
$ws.Range("E5").Value = "ay7l"
$ws.Range("F5").Value = "0tb7"
$ws.Range("G5").Value = "Changed Label"
$ws.Range("H5").Value = "👍"
$ws.Range("I5").Value = "Orig desc"
$ws.Range("J5").Value = "BOOL"
$ws.Range("K5").Value = "COUNT"

# ---------------------------------------------------------------------------
# Sheet "Entry Base" - refresh row 2/3 metadata, mark row 2 deleted with a
# note, shift row 3's period, and append row 4 as the new entry.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Entry Base")

$ws.Range("A2").Value = "lgs4l68g-0f7a"
$ws.Range("B2").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C2").Value = "lgs4l690"
$ws.Range("D2").Value = $true
$ws.Range("G2").Value = "2023-04-22T06"
$ws.Range("H2").Value = "Orig note"

$ws.Range("A3").Value = "lgs4l68h-13pq"
$ws.Range("B3").Value = "2023-04-22T15:18:42.641Z"
$ws.Range("C3").Value = "lgs4l68h"
$ws.Range("F3").Value = "lgs4l68s-gttg"
$ws.Range("G3").Value = "2023-04-22T10:18:42"

$ws.Range("A4").Value = "lgs4l690-uvob"
$ws.Range("B4").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C4").Value = "lgs4l690"
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = "ay7l"
$ws.Range("F4").Value = "lgricx7k-08al"
$ws.Range("G4").Value = "2023-04-22T06"
$ws.Range("H4").Value = "Updated noted"

# ---------------------------------------------------------------------------
# Sheet "Entry Points" - refresh row 2/3 metadata, mark row 3 deleted, and
# append row 4 as the new point value.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Entry Points")

$ws.Range("A2").Value = "lgs4l68h-w50n"
$ws.Range("B2").Value = "2023-04-22T15:18:42.641Z"
$ws.Range("C2").Value = "lgs4l68h"

$ws.Range("A3").Value = "lgs4l68h-1bns"
$ws.Range("B3").Value = "2023-04-22T15:18:42.641Z"
$ws.Range("C3").Value = "lgs4l9m5"
$ws.Range("D3").Value = $true

$ws.Range("A4").Value = "lgs4l690-06yr"
$ws.Range("B4").Value = "2023-04-22T15:18:42.641Z"
$ws.Range("C4").Value = "lgs4l690"
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = "ay7l"
$ws.Range("F4").Value = "0pc6"
$ws.Range("G4").Value = "lgricx7k-08al"
# leading apostrophe forces text storage, matching the sibling "5" cell
# (the sheet's numberStoredAsText ignored-error covers exactly this).
$ws.Range("H4").Value = "'6"
